$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New timekeeping rows (12-15), mirroring the existing row layout:
#   A = Date, B = Start Time, C = End Time, D = Interrupt (minutes),
#   E = Delta (shared formula C-B-TIME(0,D,0)), F = Activity (shared string)
# ---------------------------------------------------------------------------

$newRows = @(
    @{ Row = 12; Date = 44159; Start = 0.53819444444444442; End = 0.6166666666666667;  Interrupt = 0 },
    @{ Row = 13; Date = 44161; Start = 0.58750000000000002; End = 0.71180555555555547; Interrupt = 40 },
    @{ Row = 14; Date = 44162; Start = 0.54583333333333328; End = 0.6333333333333333;  Interrupt = 2 },
    @{ Row = 15; Date = 44163; Start = 0.76944444444444438; End = 0.78888888888888886; Interrupt = 0 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Range("A$row").Value = $r.Date
    $ws.Range("A2").Copy()
    $ws.Range("A$row").PasteSpecial(-4122) | Out-Null

    $ws.Range("B$row").Value = $r.Start
    $ws.Range("B2").Copy()
    $ws.Range("B$row").PasteSpecial(-4122) | Out-Null

    $ws.Range("C$row").Value = $r.End
    $ws.Range("C2").Copy()
    $ws.Range("C$row").PasteSpecial(-4122) | Out-Null

    $ws.Range("D$row").Value = $r.Interrupt

    $ws.Range("F$row").Value2 = $ws.Range("F2").Value2
}

# Fill column E (Delta) for the new block in one shot so the engine keeps
# it as a proper shared formula (same relative formula Excel's fill-handle
# would have produced when dragging E2's formula down through row 15).
$ws.Range("E12:E15").Formula = "=C12-B12-TIME(0,D12,0)"

# ---------------------------------------------------------------------------
# Three new formatted (but empty) cells in column G, using a new 0.00E+00
# number format (creates a new cellXfs entry).
# ---------------------------------------------------------------------------
$ws.Range("G18").NumberFormat = "0.00E+00"
$ws.Range("G19").NumberFormat = "0.00E+00"
$ws.Range("G20").NumberFormat = "0.00E+00"

# ---------------------------------------------------------------------------
# Selection / view state changes.
# ---------------------------------------------------------------------------
$ws.Range("G16:H21").Select() | Out-Null
